# Apply updated "want-to-go" counts (column F) and minimum price (column G)
# values scraped at commit 456a3b4, across the three relevant worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 642
$ws1.Range("F3").Value  = 709
$ws1.Range("F4").Value  = 946
$ws1.Range("F5").Value  = 723
$ws1.Range("F6").Value  = 839
$ws1.Range("F7").Value  = 404
$ws1.Range("F8").Value  = 605
$ws1.Range("F10").Value = 1211
$ws1.Range("F11").Value = 639
$ws1.Range("F12").Value = 384
$ws1.Range("F13").Value = 510
$ws1.Range("F14").Value = 166
$ws1.Range("F15").Value = 14
$ws1.Range("F16").Value = 508
$ws1.Range("F17").Value = 356
$ws1.Range("F18").Value = 354
$ws1.Range("G18").Value = 68
$ws1.Range("F22").Value = 579
$ws1.Range("F24").Value = 771
$ws1.Range("F25").Value = 6

# ---- Sheet "演出" ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 21

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 642
$ws4.Range("F7").Value  = 709
$ws4.Range("F8").Value  = 946
$ws4.Range("F9").Value  = 723
$ws4.Range("F10").Value = 839
$ws4.Range("F11").Value = 404
$ws4.Range("F12").Value = 605
$ws4.Range("F14").Value = 1211
$ws4.Range("F15").Value = 639
$ws4.Range("F17").Value = 21
$ws4.Range("F18").Value = 384
$ws4.Range("F19").Value = 510
$ws4.Range("F21").Value = 166
$ws4.Range("F22").Value = 14
$ws4.Range("F23").Value = 508
$ws4.Range("F25").Value = 356
$ws4.Range("F26").Value = 354
$ws4.Range("G26").Value = 68
$ws4.Range("F36").Value = 579
$ws4.Range("F38").Value = 771
$ws4.Range("F39").Value = 6
